# Update Leve profit calculation values across sheets per refreshed market data.
$wb = $excel.ActiveWorkbook

$wsALC = $wb.Worksheets.Item("ALC")
$wsALC.Cells.Item(33, 8).Value = 220.31035
$wsALC.Cells.Item(33, 9).Value = 132.88889
$wsALC.Cells.Item(33, 11).Value = 132.88889
$wsALC.Cells.Item(33, 13).Value = 96.11111
$wsALC.Cells.Item(107, 8).Value = 514.82355
$wsALC.Cells.Item(107, 9).Value = 559.5
$wsALC.Cells.Item(107, 10).Value = 407.6
$wsALC.Cells.Item(107, 11).Value = 559.5
$wsALC.Cells.Item(107, 12).Value = 407.6
$wsALC.Cells.Item(107, 13).Value = 1360.5
$wsALC.Cells.Item(107, 14).Value = -4247.6
$wsALC.Cells.Item(137, 8).Value = 1401
$wsALC.Cells.Item(137, 9).Value = 1000
$wsALC.Cells.Item(137, 10).Value = 1601.5
$wsALC.Cells.Item(137, 11).Value = 3000
$wsALC.Cells.Item(137, 12).Value = 4804.5
$wsALC.Cells.Item(137, 13).Value = -450
$wsALC.Cells.Item(137, 14).Value = -9904.5
$wsALC.Cells.Item(138, 8).Value = 3234.36
$wsALC.Cells.Item(138, 9).Value = 1132.5217
$wsALC.Cells.Item(138, 10).Value = 3862.182
$wsALC.Cells.Item(138, 11).Value = 3397.5651
$wsALC.Cells.Item(138, 12).Value = 11586.546
$wsALC.Cells.Item(138, 13).Value = 1742.4349
$wsALC.Cells.Item(138, 14).Value = -21866.546

$wsARM = $wb.Worksheets.Item("ARM")
$wsARM.Cells.Item(32, 8).Value = 8917.698
$wsARM.Cells.Item(32, 9).Value = 7630.1685
$wsARM.Cells.Item(32, 10).Value = 25287.715
$wsARM.Cells.Item(32, 11).Value = 7630.1685
$wsARM.Cells.Item(32, 12).Value = 25287.715
$wsARM.Cells.Item(32, 13).Value = -7343.1685
$wsARM.Cells.Item(32, 14).Value = -25861.715
$wsARM.Cells.Item(61, 8).Value = 1309.3793
$wsARM.Cells.Item(61, 9).Value = 1379.6923
$wsARM.Cells.Item(61, 10).Value = 700
$wsARM.Cells.Item(61, 11).Value = 1379.6923
$wsARM.Cells.Item(61, 12).Value = 700
$wsARM.Cells.Item(61, 13).Value = -1167.6923
$wsARM.Cells.Item(61, 14).Value = -1124
$wsARM.Cells.Item(74, 8).Value = 1330.069
$wsARM.Cells.Item(74, 9).Value = 1330.069
$wsARM.Cells.Item(74, 11).Value = 1330.069
$wsARM.Cells.Item(74, 13).Value = -456.069
$wsARM.Cells.Item(77, 8).Value = 1330.069
$wsARM.Cells.Item(77, 9).Value = 1330.069
$wsARM.Cells.Item(77, 11).Value = 6650.344999999999
$wsARM.Cells.Item(77, 13).Value = -2282.344999999999
$wsARM.Cells.Item(122, 8).Value = 1229.5714
$wsARM.Cells.Item(122, 9).Value = 1033.3334
$wsARM.Cells.Item(122, 10).Value = 1376.75
$wsARM.Cells.Item(122, 11).Value = 3100.0002
$wsARM.Cells.Item(122, 12).Value = 4130.25
$wsARM.Cells.Item(122, 13).Value = -650.0001999999999
$wsARM.Cells.Item(122, 14).Value = -9030.25
$wsARM.Cells.Item(136, 8).Value = 1309.3793
$wsARM.Cells.Item(136, 9).Value = 1379.6923
$wsARM.Cells.Item(136, 10).Value = 700
$wsARM.Cells.Item(136, 11).Value = 4139.0769
$wsARM.Cells.Item(136, 12).Value = 2100
$wsARM.Cells.Item(136, 13).Value = -1589.0769
$wsARM.Cells.Item(136, 14).Value = -7200
$wsARM.Cells.Item(137, 8).Value = 44900
$wsARM.Cells.Item(137, 10).Value = 44900
$wsARM.Cells.Item(137, 12).Value = 44900
$wsARM.Cells.Item(137, 14).Value = -55100

$wsCRP = $wb.Worksheets.Item("CRP")
$wsCRP.Cells.Item(31, 8).Value = 2587.2964
$wsCRP.Cells.Item(31, 9).Value = 2349.8918
$wsCRP.Cells.Item(31, 10).Value = 3104
$wsCRP.Cells.Item(31, 11).Value = 2349.8918
$wsCRP.Cells.Item(31, 12).Value = 3104
$wsCRP.Cells.Item(31, 13).Value = -2054.8918
$wsCRP.Cells.Item(31, 14).Value = -3694
$wsCRP.Cells.Item(34, 8).Value = 2587.2964
$wsCRP.Cells.Item(34, 9).Value = 2349.8918
$wsCRP.Cells.Item(34, 10).Value = 3104
$wsCRP.Cells.Item(34, 11).Value = 2349.8918
$wsCRP.Cells.Item(34, 12).Value = 3104
$wsCRP.Cells.Item(34, 13).Value = -2147.8918
$wsCRP.Cells.Item(34, 14).Value = -3508

$wsCUL = $wb.Worksheets.Item("CUL")
$wsCUL.Cells.Item(4, 8).Value = 2029
$wsCUL.Cells.Item(4, 9).Value = 145
$wsCUL.Cells.Item(4, 10).Value = 2500
$wsCUL.Cells.Item(4, 11).Value = 435
$wsCUL.Cells.Item(4, 12).Value = 7500
$wsCUL.Cells.Item(4, 13).Value = -323
$wsCUL.Cells.Item(4, 14).Value = -7724
$wsCUL.Cells.Item(12, 8).Value = 27.571428
$wsCUL.Cells.Item(12, 9).Value = 34
$wsCUL.Cells.Item(12, 10).Value = 25.818182
$wsCUL.Cells.Item(12, 11).Value = 102
$wsCUL.Cells.Item(12, 12).Value = 77.45454599999999
$wsCUL.Cells.Item(12, 13).Value = 71
$wsCUL.Cells.Item(12, 14).Value = -423.454546
$wsCUL.Cells.Item(139, 8).Value = 64888.125
$wsCUL.Cells.Item(139, 9).Value = 72729.28999999999
$wsCUL.Cells.Item(139, 10).Value = 10000
$wsCUL.Cells.Item(139, 11).Value = 218187.87
$wsCUL.Cells.Item(139, 12).Value = 30000
$wsCUL.Cells.Item(139, 13).Value = -213047.87
$wsCUL.Cells.Item(139, 14).Value = -40280

$wsGSM = $wb.Worksheets.Item("GSM")
$wsGSM.Cells.Item(102, 8).Value = 2925
$wsGSM.Cells.Item(102, 9).Value = 4000
$wsGSM.Cells.Item(102, 10).Value = 2566.6667
$wsGSM.Cells.Item(102, 11).Value = 4000
$wsGSM.Cells.Item(102, 12).Value = 2566.6667
$wsGSM.Cells.Item(102, 13).Value = -2378
$wsGSM.Cells.Item(102, 14).Value = -5810.6667
$wsGSM.Cells.Item(113, 8).Value = 4438.9414
$wsGSM.Cells.Item(113, 9).Value = 4461.8335
$wsGSM.Cells.Item(113, 10).Value = 4384
$wsGSM.Cells.Item(113, 11).Value = 4461.8335
$wsGSM.Cells.Item(113, 12).Value = 4384
$wsGSM.Cells.Item(113, 13).Value = -2291.8335
$wsGSM.Cells.Item(113, 14).Value = -8724
$wsGSM.Cells.Item(126, 8).Value = 2873.818
$wsGSM.Cells.Item(126, 9).Value = 3068
$wsGSM.Cells.Item(126, 11).Value = 9204
$wsGSM.Cells.Item(126, 13).Value = -6734

$wsLTW = $wb.Worksheets.Item("LTW")
$wsLTW.Cells.Item(22, 8).Value = 875.28
$wsLTW.Cells.Item(22, 9).Value = 304.875
$wsLTW.Cells.Item(22, 10).Value = 1143.7059
$wsLTW.Cells.Item(22, 11).Value = 304.875
$wsLTW.Cells.Item(22, 12).Value = 1143.7059
$wsLTW.Cells.Item(22, 13).Value = -9.875
$wsLTW.Cells.Item(22, 14).Value = -1733.7059
$wsLTW.Cells.Item(27, 8).Value = 875.28
$wsLTW.Cells.Item(27, 9).Value = 304.875
$wsLTW.Cells.Item(27, 10).Value = 1143.7059
$wsLTW.Cells.Item(27, 11).Value = 304.875
$wsLTW.Cells.Item(27, 12).Value = 1143.7059
$wsLTW.Cells.Item(27, 13).Value = -197.875
$wsLTW.Cells.Item(27, 14).Value = -1357.7059
$wsLTW.Cells.Item(46, 8).Value = 2667
$wsLTW.Cells.Item(46, 9).Value = 2750
$wsLTW.Cells.Item(46, 11).Value = 2750
$wsLTW.Cells.Item(46, 13).Value = -2562
$wsLTW.Cells.Item(55, 8).Value = 225.11111
$wsLTW.Cells.Item(55, 9).Value = 277.5
$wsLTW.Cells.Item(55, 10).Value = 159.625
$wsLTW.Cells.Item(55, 11).Value = 277.5
$wsLTW.Cells.Item(55, 12).Value = 159.625
$wsLTW.Cells.Item(55, 13).Value = -104.5
$wsLTW.Cells.Item(55, 14).Value = -505.625

$wsWVR = $wb.Worksheets.Item("WVR")
$wsWVR.Cells.Item(107, 8).Value = 1581.75
$wsWVR.Cells.Item(107, 9).Value = 1584.45
$wsWVR.Cells.Item(107, 10).Value = 1568.25
$wsWVR.Cells.Item(107, 11).Value = 4753.35
$wsWVR.Cells.Item(107, 12).Value = 4704.75
$wsWVR.Cells.Item(107, 13).Value = -2833.35
$wsWVR.Cells.Item(107, 14).Value = -8544.75
$wsWVR.Cells.Item(122, 8).Value = 1375
$wsWVR.Cells.Item(122, 9).Value = 1400
$wsWVR.Cells.Item(122, 10).Value = 1300
$wsWVR.Cells.Item(122, 11).Value = 4200
$wsWVR.Cells.Item(122, 12).Value = 3900
$wsWVR.Cells.Item(122, 13).Value = -1750
$wsWVR.Cells.Item(122, 14).Value = -8800
$wsWVR.Cells.Item(126, 8).Value = 1007.6842
$wsWVR.Cells.Item(126, 9).Value = 997.06665
$wsWVR.Cells.Item(126, 10).Value = 1047.5
$wsWVR.Cells.Item(126, 11).Value = 2991.19995
$wsWVR.Cells.Item(126, 12).Value = 3142.5
$wsWVR.Cells.Item(126, 13).Value = -521.1999500000002
$wsWVR.Cells.Item(126, 14).Value = -8082.5
